$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 171, shifting existing rows 171-184 down to 172-185.
$ws.Rows.Item(171).Insert()

# Populate the new row 171 with the new weekly record (same static columns as
# every other row in the sheet; D/J/K/L/M/O/P carry the new observation).
$ws.Range("A171").Value = 5
$ws.Range("B171").Value = "Macroferia Regional de Talca"
$ws.Range("C171").Value = "Maule"
$ws.Range("D171").Value = "2021-09-22"
$ws.Range("D171").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E171").Value = 7
$ws.Range("F171").Value = 100114013
$ws.Range("G171").Value = "Zanahoria"
$ws.Range("H171").Value = "Sin especificar"
$ws.Range("I171").Value = "Primera"
$ws.Range("J171").Value = 300
$ws.Range("K171").Value = 6000
$ws.Range("L171").Value = 6000
$ws.Range("M171").Value = 6000
$ws.Range("N171").Value = "$/saco 20 kilos"
$ws.Range("O171").Value = "Región de Ñuble"
$ws.Range("P171").Value = 300
$ws.Range("Q171").Value = 20
$ws.Range("R171").Value = "Hortaliza"
